# Auto-generated update script for resum_diari_meteocat.xlsx (Dades_Meteo sheet)
# Commit: Update automatic: dades i banners [2026-02-27 20:50]

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-27 20:48:16'
$r = $ws.Range('H2')
$r.Formula = '="54%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('O2').Value = '5.6 °C'
$ws.Range('E3').Value = '2026-02-27 20:48:19'
$r = $ws.Range('H3')
$r.Formula = '="40%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('N3').Value = '0.0 °C 20:13 TU'
$ws.Range('O3').Value = '4.5 °C'
$ws.Range('E4').Value = '2026-02-27 20:48:21'
$ws.Range('J4').Value = '1024.4 hPa'
$ws.Range('O4').Value = '9.7 °C'
$ws.Range('E5').Value = '2026-02-27 20:48:24'
$r = $ws.Range('H5')
$r.Formula = '="42%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('O5').Value = '4.9 °C'
$ws.Range('E6').Value = '2026-02-27 20:48:26'
$ws.Range('O6').Value = '11.0 °C'
$ws.Range('E7').Value = '2026-02-27 20:48:29'
$ws.Range('K7').Value = '13.8 MJ/m2'
$ws.Range('O7').Value = '11.6 °C'
$ws.Range('E8').Value = '2026-02-27 20:48:31'
$r = $ws.Range('H8')
$r.Formula = '="65%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('N8').Value = '8.1 °C 20:16 TU'
$ws.Range('O8').Value = '11.8 °C'
$ws.Range('E9').Value = '2026-02-27 20:48:34'
$ws.Range('E10').Value = '2026-02-27 20:48:36'
$ws.Range('E11').Value = '2026-02-27 20:48:38'
$r = $ws.Range('H11')
$r.Formula = '="71%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('O11').Value = '8.7 °C'
$ws.Range('E12').Value = '2026-02-27 20:48:41'
$ws.Range('E13').Value = '2026-02-27 20:48:43'
$ws.Range('J13').Value = '1025.4 hPa'
$ws.Range('E14').Value = '2026-02-27 20:48:46'
$r = $ws.Range('H14')
$r.Formula = '="92%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('O14').Value = '10.7 °C'
$ws.Range('E15').Value = '2026-02-27 20:48:48'
$ws.Range('E16').Value = '2026-02-27 20:48:50'
$ws.Range('N16').Value = '0.1 °C 20:24 TU'
$ws.Range('O16').Value = '2.7 °C'
$ws.Range('E17').Value = '2026-02-27 20:48:53'
$ws.Range('N17').Value = '5.0 °C 20:20 TU'
$ws.Range('O17').Value = '7.6 °C'
$ws.Range('E18').Value = '2026-02-27 20:48:55'
$ws.Range('O18').Value = '12.0 °C'
$ws.Range('E19').Value = '2026-02-27 20:48:58'
$ws.Range('O19').Value = '10.4 °C'
$ws.Range('E20').Value = '2026-02-27 20:49:00'
$r = $ws.Range('H20')
$r.Formula = '="54%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('N20').Value = '-1.2 °C 20:29 TU'
$ws.Range('O20').Value = '3.2 °C'
$ws.Range('E21').Value = '2026-02-27 20:49:03'
$ws.Range('J21').Value = '1024.1 hPa'
$ws.Range('E22').Value = '2026-02-27 20:49:05'
$r = $ws.Range('H22')
$r.Formula = '="49%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('E23').Value = '2026-02-27 20:49:08'
$ws.Range('K23').Value = '18.2 MJ/m2'
$ws.Range('N23').Value = '1.4 °C 20:07 TU'
$ws.Range('O23').Value = '3.7 °C'
$ws.Range('E24').Value = '2026-02-27 20:49:10'
$ws.Range('J24').Value = '1023.4 hPa'
$ws.Range('E25').Value = '2026-02-27 20:49:13'
$r = $ws.Range('H25')
$r.Formula = '="34%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('N25').Value = '3.0 °C 20:22 TU'
$ws.Range('O25').Value = '6.1 °C'
$ws.Range('E26').Value = '2026-02-27 20:49:15'
$ws.Range('O26').Value = '10.3 °C'
$ws.Range('E27').Value = '2026-02-27 20:49:18'
$r = $ws.Range('H27')
$r.Formula = '="42%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('K27').Value = '16.3 MJ/m2'
$ws.Range('N27').Value = '2.6 °C 20:01 TU'
$ws.Range('O27').Value = '5.6 °C'
$ws.Range('E28').Value = '2026-02-27 20:49:20'
$ws.Range('J28').Value = '1024.5 hPa'
$ws.Range('L28').Value = '14.8 km/h - 63º 20:18 TU'
$ws.Range('O28').Value = '8.1 °C'
$ws.Range('E29').Value = '2026-02-27 20:49:23'
$ws.Range('E30').Value = '2026-02-27 20:49:25'
$ws.Range('J30').Value = '1024.4 hPa'
$ws.Range('E31').Value = '2026-02-27 20:49:28'
$ws.Range('E32').Value = '2026-02-27 20:49:30'
$r = $ws.Range('H32')
$r.Formula = '="57%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('O32').Value = '7.9 °C'
$ws.Range('E33').Value = '2026-02-27 20:49:33'
$ws.Range('J33').Value = '1023.6 hPa'
$ws.Range('E34').Value = '2026-02-27 20:49:35'
$ws.Range('E35').Value = '2026-02-27 20:49:38'
$ws.Range('E36').Value = '2026-02-27 20:49:40'
$ws.Range('E37').Value = '2026-02-27 20:49:43'
$ws.Range('E38').Value = '2026-02-27 20:49:45'
$r = $ws.Range('H38')
$r.Formula = '="88%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('O38').Value = '10.2 °C'
$ws.Range('E39').Value = '2026-02-27 20:49:48'
$ws.Range('L39').Value = '40.7 km/h - 283º 20:20 TU'
$ws.Range('O39').Value = '4.6 °C'
$ws.Range('E40').Value = '2026-02-27 20:49:50'
$ws.Range('E41').Value = '2026-02-27 20:49:52'
$r = $ws.Range('H41')
$r.Formula = '="85%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('E42').Value = '2026-02-27 20:49:55'
$ws.Range('E43').Value = '2026-02-27 20:49:57'
$r = $ws.Range('H43')
$r.Formula = '="73%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('E44').Value = '2026-02-27 20:49:59'
$ws.Range('O44').Value = '2.0 °C'
$ws.Range('E45').Value = '2026-02-27 20:50:02'
$r = $ws.Range('H45')
$r.Formula = '="43%"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range('K45').Value = '16.8 MJ/m2'
$ws.Range('E46').Value = '2026-02-27 20:50:04'
